$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.683.90"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.053.50"
$ws.Range("E3").Value = "  +4.48%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "201.98"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "632.77"
$ws.Range("E6").Value = "  +5.60%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.554"
$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.204"
$ws.Range("E9").Value = "  +3.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.053.10"
$ws.Range("E10").Value = "  +4.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.436"
$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.09"
$ws.Range("E13").Value = "  +4.17%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.609.41"
$ws.Range("E14").Value = "  +4.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.56"
$ws.Range("E15").Value = "  +7.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.519.41"
$ws.Range("E16").Value = "  +0.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000190"
$ws.Range("E17").Value = "  -0.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.027.45"
$ws.Range("E18").Value = "  +3.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.65"
$ws.Range("E19").Value = "  +6.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.01"
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.33"
$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.38"
$ws.Range("E22").Value = "  +4.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.29"
$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.42"
$ws.Range("E24").Value = "  +2.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.154.24"
$ws.Range("E25").Value = "  +2.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.42"
$ws.Range("E26").Value = "  +4.85%  "

$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.93"
$ws.Range("E28").Value = "  +2.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000110"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.40"
$ws.Range("E31").Value = "  +8.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("E32").Value = "  -0.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "512.95"
$ws.Range("E33").Value = "  +0.97%  "

$ws.Range("E34").Value = "  +7.61%  "

$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "21.08"
$ws.Range("E36").Value = "  +4.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.39"
$ws.Range("E37").Value = "  -0.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.384"
$ws.Range("E38").Value = "  +11.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.01"
$ws.Range("E39").Value = "  +1.73%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "189.34"
$ws.Range("E40").Value = "  +4.62%  "

$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.105"
$ws.Range("E41").Value = "  -4.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.113"
$ws.Range("E42").Value = "  -0.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.07"
$ws.Range("E44").Value = "  +1.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.06"
$ws.Range("E45").Value = "  +7.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.26"
$ws.Range("E46").Value = "  +4.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.67"
$ws.Range("E47").Value = "  +0.45%  "

$ws.Range("E48").Value = "  +7.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.709"
$ws.Range("E49").Value = "  +6.84%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.35"
$ws.Range("E50").Value = "  +0.70%  "

$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.90"
$ws.Range("E51").Value = "  +4.86%  "
